$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.33650637266031
$ws.Range("C2").Value = 4.408983019100955
$ws.Range("E2").Value = 9.332530527201342
$ws.Range("F2").Value = 65.41654817454751
$ws.Range("G2").Value = 3.826983990749906
$ws.Range("J2").Value = 13.71773040068543
$ws.Range("K2").Value = 16.2391337300935
$ws.Range("L2").Value = 9.428076638948738
$ws.Range("M2").Value = 16.97206394233627
$ws.Range("B3").Value = 19.3573741040462
$ws.Range("C3").Value = 4.178465634580479
$ws.Range("E3").Value = 9.358939755907141
$ws.Range("F3").Value = 64.50928574972426
$ws.Range("G3").Value = 3.830557660038769
$ws.Range("J3").Value = 13.60354891776379
$ws.Range("K3").Value = 16.24324354202467
$ws.Range("L3").Value = 9.497999206220701
$ws.Range("M3").Value = 17.07228777481874
$ws.Range("B4").Value = 19.37656023124551
$ws.Range("C4").Value = 4.029499198021233
$ws.Range("E4").Value = 9.377228553769044
$ws.Range("F4").Value = 63.94802761159888
$ws.Range("G4").Value = 3.832864079511856
$ws.Range("J4").Value = 13.53190880769038
$ws.Range("K4").Value = 16.25054958353237
$ws.Range("L4").Value = 9.543212465027326
$ws.Range("M4").Value = 17.13793353153811
$ws.Range("B5").Value = 19.38597833338346
$ws.Range("C5").Value = 3.966955155171505
$ws.Range("E5").Value = 9.38520125658615
$ws.Range("F5").Value = 63.71842436778699
$ws.Range("G5").Value = 3.833832280737706
$ws.Range("J5").Value = 13.50233630701626
$ws.Range("K5").Value = 16.25472863469352
$ws.Range("L5").Value = 9.56221258751677
$ws.Range("M5").Value = 17.1657185806406
$ws.Range("B6").Value = 19.38763870729723
$ws.Range("C6").Value = 3.956459655996208
$ws.Range("E6").Value = 9.386556471324143
$ws.Range("F6").Value = 63.68025036775087
$ws.Range("G6").Value = 3.833994763292113
$ws.Range("J6").Value = 13.49740314710213
$ws.Range("K6").Value = 16.25549510754839
$ws.Range("L6").Value = 9.565402349159266
$ws.Range("M6").Value = 17.17039473755572
$ws.Range("B7").Value = 19.37668077511477
$ws.Range("C7").Value = 4.028663110685356
$ws.Range("E7").Value = 9.377333973563386
$ws.Range("F7").Value = 63.94493446566381
$ws.Range("G7").Value = 3.832877022215549
$ws.Range("J7").Value = 13.5315115064984
$ws.Range("K7").Value = 16.25060107947319
$ws.Range("L7").Value = 9.543466375263344
$ws.Range("M7").Value = 17.13830406316345
$ws.Range("B8").Value = 19.34237727322148
$ws.Range("C8").Value = 4.331057697619971
$ws.Range("E8").Value = 9.341205198754691
$ws.Range("F8").Value = 65.10469577611651
$ws.Range("G8").Value = 3.828192976965495
$ws.Range("J8").Value = 13.67868264810732
$ws.Range("K8").Value = 16.23955788651893
$ws.Range("L8").Value = 9.45171375977371
$ws.Range("M8").Value = 17.00576905216099
$ws.Range("B9").Value = 19.32576938470683
$ws.Range("C9").Value = 4.864138021041764
$ws.Range("E9").Value = 9.286874021962079
$ws.Range("F9").Value = 67.33837448171263
$ws.Range("G9").Value = 3.819892591920814
$ws.Range("J9").Value = 13.95483956806828
$ws.Range("K9").Value = 16.25585459204096
$ws.Range("L9").Value = 9.289795958502276
$ws.Range("M9").Value = 16.77842377052178
$ws.Range("B10").Value = 19.34451558542638
$ws.Range("C10").Value = 5.218288098301433
$ws.Range("E10").Value = 9.257112915186685
$ws.Range("F10").Value = 68.94510184633737
$ws.Range("G10").Value = 3.814326791688248
$ws.Range("J10").Value = 14.14977123769422
$ws.Range("K10").Value = 16.29092241935172
$ws.Range("L10").Value = 9.181693799982103
$ws.Range("M10").Value = 16.63118844913319
$ws.Range("B11").Value = 19.35975551996861
$ws.Range("C11").Value = 5.371125168952651
$ws.Range("E11").Value = 9.245796319971408
$ws.Range("F11").Value = 69.66661905827898
$ws.Range("G11").Value = 3.811908886029353
$ws.Range("J11").Value = 14.23664416414451
$ws.Range("K11").Value = 16.31186423048198
$ws.Range("L11").Value = 9.134847940157179
$ws.Range("M11").Value = 16.56849606629667
$ws.Range("B12").Value = 19.36648914791748
$ws.Range("C12").Value = 5.427804264929991
$ws.Range("E12").Value = 9.241831810934897
$ws.Range("F12").Value = 69.93832943772978
$ws.Range("G12").Value = 3.811009566069229
$ws.Range("J12").Value = 14.26927479388019
$ws.Range("K12").Value = 16.32050814877449
$ws.Range("L12").Value = 9.117441829459942
$ws.Range("M12").Value = 16.54537183496704
$ws.Range("B13").Value = 19.36499618233032
$ws.Range("C13").Value = 5.415650748462767
$ws.Range("E13").Value = 9.24267134904186
$ws.Range("F13").Value = 69.87988153776979
$ws.Range("G13").Value = 3.811202527953878
$ws.Range("J13").Value = 14.26225915741309
$ws.Range("K13").Value = 16.3186148525988
$ws.Range("L13").Value = 9.121175745621265
$ws.Range("M13").Value = 16.55032465798398
$ws.Range("B14").Value = 19.3602902266504
$ws.Range("C14").Value = 5.375812219908281
$ws.Range("E14").Value = 9.245463719956373
$ws.Range("F14").Value = 69.68900405678384
$ws.Range("G14").Value = 3.811834572592676
$ws.Range("J14").Value = 14.2393340820758
$ws.Range("K14").Value = 16.31256109153636
$ws.Range("L14").Value = 9.133409256792367
$ws.Range("M14").Value = 16.56658127730266
$ws.Range("B15").Value = 19.35753295481717
$ws.Range("C15").Value = 5.351253886964296
$ws.Range("E15").Value = 9.247215948853652
$ws.Range("F15").Value = 69.57188439588658
$ws.Range("G15").Value = 3.812223836298977
$ws.Range("J15").Value = 14.22525692837018
$ws.Range("K15").Value = 16.30894580834737
$ws.Range("L15").Value = 9.140946003999538
$ws.Range("M15").Value = 16.57661914254351
$ws.Range("B16").Value = 19.34365448866038
$ws.Range("C16").Value = 5.208132723016313
$ws.Range("E16").Value = 9.257897307311316
$ws.Range("F16").Value = 68.89774604497585
$ws.Range("G16").Value = 3.81448709256671
$ws.Range("J16").Value = 14.14405701063129
$ws.Range("K16").Value = 16.28965391937374
$ws.Range("L16").Value = 9.184802033685258
$ws.Range("M16").Value = 16.6353717587154
$ws.Range("B17").Value = 19.33685814910445
$ws.Range("C17").Value = 5.118208451801549
$ws.Range("E17").Value = 9.265020028298117
$ws.Range("F17").Value = 68.48166166372506
$ws.Range("G17").Value = 3.815904651294297
$ws.Range("J17").Value = 14.09377677283964
$ws.Range("K17").Value = 16.27909445345112
$ws.Range("L17").Value = 9.21230192784383
$ws.Range("M17").Value = 16.67251204481573
$ws.Range("B18").Value = 19.33358108477447
$ws.Range("C18").Value = 5.065708523736382
$ws.Range("E18").Value = 9.269325888251217
$ws.Range("F18").Value = 68.24146636402465
$ws.Range("G18").Value = 3.816730730527688
$ws.Range("J18").Value = 14.06468766976589
$ws.Range("K18").Value = 16.2734908214837
$ws.Range("L18").Value = 9.228338583903763
$ws.Range("M18").Value = 16.69427755714581
$ws.Range("B19").Value = 19.33258013671988
$ws.Range("C19").Value = 5.047799637879011
$ws.Range("E19").Value = 9.270819642500081
$ws.Range("F19").Value = 68.15999518265549
$ws.Range("G19").Value = 3.817012274056618
$ws.Range("J19").Value = 14.05480974930772
$ws.Range("K19").Value = 16.27167433261452
$ws.Range("L19").Value = 9.233806062102829
$ws.Range("M19").Value = 16.70171627956365
$ws.Range("B20").Value = 19.3375162359365
$ws.Range("C20").Value = 5.127861616310996
$ws.Range("E20").Value = 9.264240155777083
$ws.Range("F20").Value = 68.52604626275016
$ws.Range("G20").Value = 3.815752639283223
$ws.Range("J20").Value = 14.09914675710533
$ws.Range("K20").Value = 16.28016991896389
$ws.Range("L20").Value = 9.209351818029431
$ws.Range("M20").Value = 16.66851664970341
$ws.Range("B21").Value = 19.36164638164907
$ws.Range("C21").Value = 5.387546300474619
$ws.Range("E21").Value = 9.244634815184563
$ws.Range("F21").Value = 69.74511171375514
$ws.Range("G21").Value = 3.811648484534819
$ws.Range("J21").Value = 14.24607501448837
$ws.Range("K21").Value = 16.31431989211948
$ws.Range("L21").Value = 9.1298069430848
$ws.Range("M21").Value = 16.56178959745835
$ws.Range("B22").Value = 19.38302557506718
$ws.Range("C22").Value = 5.550286100500881
$ws.Range("E22").Value = 9.233692253776923
$ws.Range("F22").Value = 70.53295411558784
$ws.Range("G22").Value = 3.809061075006422
$ws.Range("J22").Value = 14.34054594726206
$ws.Range("K22").Value = 16.3407965778882
$ws.Range("L22").Value = 9.079762305621436
$ws.Range("M22").Value = 16.49562782666005
$ws.Range("B23").Value = 19.37110304555981
$ws.Range("C23").Value = 5.464069497961826
$ws.Range("E23").Value = 9.239360918981758
$ws.Range("F23").Value = 70.11333308220399
$ws.Range("G23").Value = 3.810433376169184
$ws.Range("J23").Value = 14.29026963748117
$ws.Range("K23").Value = 16.32628649918225
$ws.Range("L23").Value = 9.106294884164782
$ws.Range("M23").Value = 16.5306111437283
$ws.Range("B24").Value = 19.33721675153757
$ws.Range("C24").Value = 5.123499915825992
$ws.Range("E24").Value = 9.264592079269244
$ws.Range("F24").Value = 68.50598303441818
$ws.Range("G24").Value = 3.815821329311035
$ws.Range("J24").Value = 14.09671955465779
$ws.Range("K24").Value = 16.27968224560843
$ws.Range("L24").Value = 9.21068485668426
$ws.Range("M24").Value = 16.67032168111406
$ws.Range("B25").Value = 19.32482790083174
$ws.Range("C25").Value = 4.726456494730991
$ws.Range("E25").Value = 9.299794540713505
$ws.Range("F25").Value = 66.73945639262186
$ws.Range("G25").Value = 3.822044047863149
$ws.Range("J25").Value = 13.88150291403199
$ws.Range("K25").Value = 16.24738387718471
$ws.Range("L25").Value = 9.331683514502627
$ws.Range("M25").Value = 16.83644697949977
